# Auto-generated edit script applying scheduled-runner profit recalculations
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4545.875
$ws.Range("I33").Value = 6657.875
$ws.Range("J33").Value = 321.875
$ws.Range("K33").Value = 6657.875
$ws.Range("L33").Value = 321.875
$ws.Range("M33").Value = -6428.875
$ws.Range("N33").Value = -779.875
$ws.Range("H64").Value = 3750
$ws.Range("H67").Value = 3750
$ws.Range("H100").Value = 334517.34
$ws.Range("I100").Value = 500901
$ws.Range("K100").Value = 500901
$ws.Range("M100").Value = -500360
$ws.Range("H103").Value = 632.375
$ws.Range("I103").Value = 1239
$ws.Range("J103").Value = 545.7143
$ws.Range("K103").Value = 3717
$ws.Range("L103").Value = 1637.1429
$ws.Range("M103").Value = -3131
$ws.Range("N103").Value = -2809.1429
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = $null
$ws.Range("H116").Value = 10082.2
$ws.Range("J116").Value = 11537.444
$ws.Range("L116").Value = 11537.444
$ws.Range("N116").Value = -18421.444
$ws.Range("H135").Value = 1085.4642
$ws.Range("I135").Value = 971.6818
$ws.Range("J135").Value = 1502.6666
$ws.Range("K135").Value = 8745.136199999999
$ws.Range("L135").Value = 13523.9994
$ws.Range("M135").Value = -6210.136199999999
$ws.Range("N135").Value = -18593.9994
$ws.Range("H138").Value = 4230.64
$ws.Range("I138").Value = 3280.5833
$ws.Range("J138").Value = 5107.615
$ws.Range("K138").Value = 9841.749899999999
$ws.Range("L138").Value = 15322.845
$ws.Range("M138").Value = -4701.749899999999
$ws.Range("N138").Value = -25602.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2448.889
$ws.Range("I2").Value = 2364.8
$ws.Range("J2").Value = 2689.1428
$ws.Range("K2").Value = 2364.8
$ws.Range("L2").Value = 2689.1428
$ws.Range("M2").Value = -2251.8
$ws.Range("N2").Value = -2915.1428
$ws.Range("H32").Value = 6070.316
$ws.Range("I32").Value = 4577.5557
$ws.Range("K32").Value = 4577.5557
$ws.Range("M32").Value = -4290.5557
$ws.Range("H102").Value = 2873.6
$ws.Range("I102").Value = 2808.6667
$ws.Range("K102").Value = 2808.6667
$ws.Range("M102").Value = -1186.6667
$ws.Range("H110").Value = 1431.5555
$ws.Range("I110").Value = 1449.4
$ws.Range("J110").Value = 1342.3334
$ws.Range("K110").Value = 1449.4
$ws.Range("L110").Value = 1342.3334
$ws.Range("M110").Value = 595.5999999999999
$ws.Range("N110").Value = -5432.3334
$ws.Range("H116").Value = 2448.889
$ws.Range("I116").Value = 2364.8
$ws.Range("J116").Value = 2689.1428
$ws.Range("K116").Value = 2364.8
$ws.Range("L116").Value = 2689.1428
$ws.Range("M116").Value = -70.80000000000018
$ws.Range("N116").Value = -7277.1428
$ws.Range("H122").Value = 1734.4546
$ws.Range("I122").Value = 1412.875
$ws.Range("J122").Value = 2592
$ws.Range("K122").Value = 4238.625
$ws.Range("L122").Value = 7776
$ws.Range("M122").Value = -1788.625
$ws.Range("N122").Value = -12676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2448.889
$ws.Range("I3").Value = 2364.8
$ws.Range("J3").Value = 2689.1428
$ws.Range("K3").Value = 2364.8
$ws.Range("L3").Value = 2689.1428
$ws.Range("M3").Value = -2250.8
$ws.Range("N3").Value = -2917.1428
$ws.Range("H81").Value = 20849
$ws.Range("J81").Value = 20849
$ws.Range("L81").Value = 20849
$ws.Range("N81").Value = -22971
$ws.Range("H84").Value = 20849
$ws.Range("J84").Value = 20849
$ws.Range("L84").Value = 62547
$ws.Range("N84").Value = -73155
$ws.Range("H99").Value = 2261.818
$ws.Range("I99").Value = 2153.3333
$ws.Range("K99").Value = 2153.3333
$ws.Range("M99").Value = -655.3332999999998
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null
$ws.Range("H107").Value = 1329.0454
$ws.Range("I107").Value = 1316.8125
$ws.Range("J107").Value = 1361.6666
$ws.Range("K107").Value = 1316.8125
$ws.Range("L107").Value = 1361.6666
$ws.Range("M107").Value = 603.1875
$ws.Range("N107").Value = -5201.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 29590.166
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31250
$ws.Range("H58").Value = 2450.4
$ws.Range("J58").Value = 2721.5715
$ws.Range("L58").Value = 2721.5715
$ws.Range("N58").Value = -3127.5715
$ws.Range("H62").Value = 501503
$ws.Range("I62").Value = 1000000
$ws.Range("J62").Value = 3006
$ws.Range("K62").Value = 1000000
$ws.Range("L62").Value = 3006
$ws.Range("M62").Value = -999376
$ws.Range("N62").Value = -4254
$ws.Range("H65").Value = 501503
$ws.Range("I65").Value = 1000000
$ws.Range("J65").Value = 3006
$ws.Range("K65").Value = 5000000
$ws.Range("L65").Value = 15030
$ws.Range("M65").Value = -4996880
$ws.Range("N65").Value = -21270
$ws.Range("H122").Value = 2400
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("H134").Value = 3691.0435
$ws.Range("I134").Value = 3915.5789
$ws.Range("K134").Value = 11746.7367
$ws.Range("M134").Value = -9211.736699999999
$ws.Range("H136").Value = 2450.4
$ws.Range("J136").Value = 2721.5715
$ws.Range("L136").Value = 8164.7145
$ws.Range("N136").Value = -13264.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 198.83333
$ws.Range("I97").Value = 148.5
$ws.Range("J97").Value = 299.5
$ws.Range("K97").Value = 445.5
$ws.Range("L97").Value = 898.5
$ws.Range("M97").Value = 50.5
$ws.Range("N97").Value = -1890.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11159.4
$ws.Range("I80").Value = 1199
$ws.Range("K80").Value = 1199
$ws.Range("M80").Value = -201
$ws.Range("H83").Value = 11159.4
$ws.Range("I83").Value = 1199
$ws.Range("K83").Value = 5995
$ws.Range("M83").Value = -1003
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("H122").Value = 2009.2858
$ws.Range("I122").Value = 1266.4706
$ws.Range("K122").Value = 3799.4118
$ws.Range("M122").Value = -1349.4118
$ws.Range("H126").Value = 3816.111
$ws.Range("J126").Value = 3976.923
$ws.Range("L126").Value = 11930.769
$ws.Range("N126").Value = -16870.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2750.75
$ws.Range("I40").Value = 2001
$ws.Range("K40").Value = 2001
$ws.Range("M40").Value = -1865
$ws.Range("H82").Value = 1904.2
$ws.Range("J82").Value = 2632.1667
$ws.Range("L82").Value = 2632.1667
$ws.Range("N82").Value = -3354.1667
$ws.Range("H85").Value = 1904.2
$ws.Range("J85").Value = 2632.1667
$ws.Range("L85").Value = 2632.1667
$ws.Range("N85").Value = -5128.1667
$ws.Range("H93").Value = 1851.5862
$ws.Range("I93").Value = 1894.0834
$ws.Range("K93").Value = 1894.0834
$ws.Range("M93").Value = -646.0834
$ws.Range("H122").Value = 3344.4285
$ws.Range("I122").Value = 2991.375
$ws.Range("J122").Value = 3815.1667
$ws.Range("K122").Value = 8974.125
$ws.Range("L122").Value = 11445.5001
$ws.Range("M122").Value = -6524.125
$ws.Range("N122").Value = -16345.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 25000
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("N43").Value = -25298
$ws.Range("H107").Value = 2633.4285
$ws.Range("J107").Value = 4511
$ws.Range("L107").Value = 13533
$ws.Range("N107").Value = -17373
$ws.Range("H122").Value = 6405.857
$ws.Range("I122").Value = 6547.6665
$ws.Range("K122").Value = 19642.9995
$ws.Range("M122").Value = -17192.9995
$ws.Range("H132").Value = 11012.5
$ws.Range("I132").Value = 7750
$ws.Range("J132").Value = 27325
$ws.Range("K132").Value = 23250
$ws.Range("L132").Value = 81975
$ws.Range("M132").Value = -20720
$ws.Range("N132").Value = -87035
